$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "MessageWithQuote"
$ws.Range("B4").Value = "J'aime les écoeurants ""et toi"""

$ws.Range("B5").Select()
